$d = $word.ActiveDocument

$d.Content.Find.Execute("43÷5=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "28÷3=9, 1", 2)
$d.Content.Find.Execute("22÷9=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "91÷3=30, 1", 2)
$d.Content.Find.Execute("99÷8=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=21, 3", 2)
$d.Content.Find.Execute("62÷9=6, 8", $true, $false, $false, $false, $false, $true, 1, $false, "71÷3=23, 2", 2)
$d.Content.Find.Execute("61÷9=6, 7", $true, $false, $false, $false, $false, $true, 1, $false, "10÷5=2, 0", 2)
$d.Content.Find.Execute("13÷8=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=3, 7", 2)
$d.Content.Find.Execute("79÷9=8, 7", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=5, 3", 2)
$d.Content.Find.Execute("36÷9=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "64÷7=9, 1", 2)
$d.Content.Find.Execute("21÷9=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=4, 2", 2)
$d.Content.Find.Execute("58÷9=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "91÷6=15, 1", 2)
$d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "27÷7=3, 6", 2)
$d.Content.Find.Execute("54÷4=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "22÷6=3, 4", 2)
$d.Content.Find.Execute("83÷7=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 2)
$d.Content.Find.Execute("76÷8=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "44÷6=7, 2", 2)
$d.Content.Find.Execute("32÷4=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "40÷8=5, 0", 2)
$d.Content.Find.Execute("28÷4=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=10, 1", 2)
$d.Content.Find.Execute("83÷3=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "72÷4=18, 0", 2)
$d.Content.Find.Execute("62÷3=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2)
$d.Content.Find.Execute("34÷9=3, 7", $true, $false, $false, $false, $false, $true, 1, $false, "64÷5=12, 4", 2)
$d.Content.Find.Execute("49÷9=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=18, 3", 2)
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "84÷7=12, 0", 2)
$d.Content.Find.Execute("34÷5=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "15÷4=3, 3", 2)
$d.Content.Find.Execute("42÷7=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "80÷2=40, 0", 2)
$d.Content.Find.Execute("38÷3=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "97÷7=13, 6", 2)
